# Login_TestSuite.xlsx - "updated code for dynamic xml creator and dynamic class creator"
#
# The original 3 "Login" scenarios (valid user / invalid user) are replaced
# with 3 "Registration / Login" scenarios that exercise the
# way2automationHomePage registration flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Clear out the cells that are no longer populated in the new layout
# ---------------------------------------------------------------------
$cellsToClear = @("A9", "A16", "D3", "F3", "F10", "F17")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------
# 2. Write the new cell values (header row is unchanged).
#    Each test-case block is filled in Steps/method_Name/testdata order
#    first (D, E, F), then the ScenarioName cell (C) is revisited - this
#    mirrors the order the rows were actually authored in the workbook.
# ---------------------------------------------------------------------

# --- Test case 1 (rows 2-7) ---
$ws.Range("D2").Value  = "User is on way2automationHomePage"
$ws.Range("E2").Value  = "Comman_Reusables.launchURL"
$ws.Range("E3").Value  = "Home_Page.validateUserIsOnHomePage"
$ws.Range("D4").Value  = "User clicks on the registration link"
$ws.Range("E4").Value  = "Home_Page.clickRegistrationLink"
$ws.Range("D5").Value  = "Registration pop up is displayed"
$ws.Range("E5").Value  = "Home_Page.validateRegistrationFormDisplayed"
$ws.Range("D6").Value  = "Register with a valid user"
$ws.Range("E6").Value  = "home_Page.registrationOnPopUp"
$ws.Range("F6").Value  = "Registration_testdata|Registration_Details|1"
$ws.Range("E7").Value  = "Registration_Page.registrationOnRegistrationPage"
$ws.Range("F7").Value  = "Registration_testdata|Registration_Details|1-2"
$ws.Range("C2").Value  = "Login First Test Case"

# --- Test case 2 (rows 9-14) ---
$ws.Range("D9").Value  = "User is on way2automationHomePage"
$ws.Range("E9").Value  = "Comman_Reusables.launchURL"
$ws.Range("E10").Value = "Home_Page.validateUserIsOnHomePage"
$ws.Range("D11").Value = "User clicks on the registration link"
$ws.Range("E11").Value = "Home_Page.clickRegistrationLink"
$ws.Range("D12").Value = "Registration pop up is displayed"
$ws.Range("E12").Value = "Home_Page.validateRegistrationFormDisplayed"
$ws.Range("D13").Value = "Register with a valid user"
$ws.Range("E13").Value = "home_Page.registrationOnPopUp"
$ws.Range("F13").Value = "Registration_testdata|Registration_Details|1"
$ws.Range("E14").Value = "Registration_Page.registrationOnRegistrationPage"
$ws.Range("F14").Value = "Registration_testdata|Registration_Details|1-1"
$ws.Range("C9").Value  = "Login Second Test Case"

# --- Test case 3 (rows 16-21) ---
$ws.Range("D16").Value = "User is on way2automationHomePage"
$ws.Range("E16").Value = "Comman_Reusables.launchURL"
$ws.Range("E17").Value = "Home_Page.validateUserIsOnHomePage"
$ws.Range("D18").Value = "User clicks on the registration link"
$ws.Range("E18").Value = "Home_Page.clickRegistrationLink"
$ws.Range("D19").Value = "Registration pop up is displayed"
$ws.Range("E19").Value = "Home_Page.validateRegistrationFormDisplayed"
$ws.Range("D20").Value = "Register with a valid user"
$ws.Range("E20").Value = "home_Page.registrationOnPopUp"
$ws.Range("F20").Value = "Registration_testdata|Registration_Details|1"
$ws.Range("E21").Value = "Registration_Page.registrationOnRegistrationPage"
$ws.Range("F21").Value = "Registration_testdata|Registration_Details|1-1"
$ws.Range("C16").Value = "Login Third Test Case"

# ---------------------------------------------------------------------
# 3. Widen columns E and F so the longer method/testdata strings fit
#    (mirrors the workbook author re-running best-fit autosize)
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 46.7109375
$ws.Columns.Item(6).ColumnWidth = 43.7109375

# ---------------------------------------------------------------------
# 4. Move the saved cursor/selection from C16 to A9
# ---------------------------------------------------------------------
$ws.Range("A9").Select()

$wb.Save()
